# Update vm_pu.xlsx results for Case_0_7 (380 kV case) in res_bus sheet.
# Column B (slack bus vm_pu setpoint) goes from 1.05 to 1.02 p.u.,
# and the downstream bus voltage results (C:F, I:N) are refreshed
# with the corresponding load-flow solution for the new setpoint.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040240711994879
$ws.Range("D2").Value = 1.040678017384052
$ws.Range("E2").Value = 1.047881117036417
$ws.Range("F2").Value = 1.056895564117549
$ws.Range("I2").Value = 1.036382026818991
$ws.Range("J2").Value = 1.045328385389583
$ws.Range("K2").Value = 1.043459685592757
$ws.Range("L2").Value = 1.0506425168711
$ws.Range("M2").Value = 1.059632020472387
$ws.Range("N2").Value = 1.018961036247124

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041362358261452
$ws.Range("D3").Value = 1.041657361855371
$ws.Range("E3").Value = 1.048882162365443
$ws.Range("F3").Value = 1.057991224953231
$ws.Range("I3").Value = 1.036566271504541
$ws.Range("J3").Value = 1.046094708349907
$ws.Range("K3").Value = 1.044249058424736
$ws.Range("L3").Value = 1.051455005627237
$ws.Range("M3").Value = 1.060540693633362
$ws.Range("N3").Value = 1.019220374192482

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042088158385703
$ws.Range("D4").Value = 1.042291352682734
$ws.Range("E4").Value = 1.04953026616906
$ws.Range("F4").Value = 1.058700620264352
$ws.Range("I4").Value = 1.036684005054374
$ws.Range("J4").Value = 1.046590043741819
$ws.Range("K4").Value = 1.04475951062213
$ws.Range("L4").Value = 1.051980488967751
$ws.Range("M4").Value = 1.061128501162442
$ws.Range("N4").Value = 1.019387877075953

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04239329020691
$ws.Range("D5").Value = 1.042557951687871
$ws.Range("E5").Value = 1.049802815287056
$ws.Range("F5").Value = 1.058998952935273
$ws.Range("I5").Value = 1.036733144571142
$ws.Range("J5").Value = 1.046798156697497
$ws.Range("K5").Value = 1.044974026996962
$ws.Range("L5").Value = 1.052201342048281
$ws.Range("M5").Value = 1.061375576034749
$ws.Range("N5").Value = 1.019458221941104

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042444523521239
$ws.Range("D6").Value = 1.042602718894354
$ws.Range("E6").Value = 1.049848582541693
$ws.Range("F6").Value = 1.05904905033835
$ws.Range("I6").Value = 1.036741374451764
$ws.Range("J6").Value = 1.046833092389472
$ws.Range("K6").Value = 1.045010040683059
$ws.Range("L6").Value = 1.052238420727831
$ws.Range("M6").Value = 1.061417058682768
$ws.Range("N6").Value = 1.019470028850519

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042092235549783
$ws.Range("D7").Value = 1.042294914719587
$ws.Range("E7").Value = 1.049533907644741
$ws.Range("F7").Value = 1.058704606195766
$ws.Range("I7").Value = 1.036684663056702
$ws.Range("J7").Value = 1.04659282505249
$ws.Range("K7").Value = 1.044762377305816
$ws.Range("L7").Value = 1.051983440253607
$ws.Range("M7").Value = 1.061131802743805
$ws.Range("N7").Value = 1.019388817315769

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040619773862835
$ws.Range("D8").Value = 1.041008931475077
$ws.Range("E8").Value = 1.048219350320266
$ws.Range("F8").Value = 1.057265758936826
$ws.Range("I8").Value = 1.036444600587534
$ws.Range("J8").Value = 1.045587476946346
$ws.Range("K8").Value = 1.043726525110038
$ws.Range("L8").Value = 1.050917153281593
$ws.Range("M8").Value = 1.059939144957766
$ws.Range("N8").Value = 1.01904874406362

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038025231241156
$ws.Range("D9").Value = 1.038745082754918
$ws.Range("E9").Value = 1.045905694591405
$ws.Range("F9").Value = 1.054733612202437
$ws.Range("I9").Value = 1.036010211689518
$ws.Range("J9").Value = 1.043811886558991
$ws.Range("K9").Value = 1.041898731646365
$ws.Range("L9").Value = 1.0490362937449
$ws.Range("M9").Value = 1.057836265620478
$ws.Range("N9").Value = 1.018447149742212

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036295574027511
$ws.Range("D10").Value = 1.037237339989424
$ws.Range("E10").Value = 1.044365112803725
$ws.Range("F10").Value = 1.053047715713274
$ws.Range("I10").Value = 1.035712982535733
$ws.Range("J10").Value = 1.042625436289221
$ws.Range("K10").Value = 1.040678523511222
$ws.Range("L10").Value = 1.047781087660609
$ws.Range("M10").Value = 1.056433490490162
$ws.Range("N10").Value = 1.018044514957067

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035546611484857
$ws.Range("D11").Value = 1.036584822742202
$ws.Range("E11").Value = 1.043698462569463
$ws.Range("F11").Value = 1.052318223256376
$ws.Range("I11").Value = 1.035582468801543
$ws.Range("J11").Value = 1.042111041664611
$ws.Range("K11").Value = 1.040149759615425
$ws.Range("L11").Value = 1.047237259318742
$ws.Range("M11").Value = 1.055825868226475
$ws.Range("N11").Value = 1.017869797050751

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035268410970215
$ws.Range("D12").Value = 1.0363425004418
$ws.Range("E12").Value = 1.043450903642998
$ws.Range("F12").Value = 1.052047333822726
$ws.Range("I12").Value = 1.035533718050626
$ws.Range("J12").Value = 1.04191987404046
$ws.Range("K12").Value = 1.039953291931597
$ws.Range("L12").Value = 1.047035209539012
$ws.Range("M12").Value = 1.055600138206805
$ws.Range("N12").Value = 1.017804842772705

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.03532808605389
$ws.Range("D13").Value = 1.03639447706374
$ws.Range("E13").Value = 1.04350400295105
$ws.Range("F13").Value = 1.052105437067646
$ws.Range("I13").Value = 1.035544187566631
$ws.Range("J13").Value = 1.04196088460995
$ws.Range("K13").Value = 1.0399954376919
$ws.Range("L13").Value = 1.047078552066406
$ws.Range("M13").Value = 1.055648559517732
$ws.Range("N13").Value = 1.017818778232444

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035523615400127
$ws.Range("D14").Value = 1.036564791252728
$ws.Range("E14").Value = 1.043677997953788
$ws.Range("F14").Value = 1.052295829891831
$ws.Range("I14").Value = 1.035578444600666
$ws.Range("J14").Value = 1.042095241701593
$ws.Range("K14").Value = 1.040133520786052
$ws.Range("L14").Value = 1.047220558795954
$ws.Range("M14").Value = 1.055807209968313
$ws.Range("N14").Value = 1.017864429057846

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035644087123989
$ws.Range("D15").Value = 1.036669734284243
$ws.Range("E15").Value = 1.0437852105779
$ws.Range("F15").Value = 1.052413147326995
$ws.Range("I15").Value = 1.035599515429689
$ws.Range("J15").Value = 1.042178010470534
$ws.Range("K15").Value = 1.040218590215914
$ws.Range("L15").Value = 1.047308047497426
$ws.Range("M15").Value = 1.055904955490262
$ws.Range("N15").Value = 1.017892548581306

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036345279914435
$ws.Range("D16").Value = 1.037280652666697
$ws.Range("E16").Value = 1.044409365263379
$ws.Range("F16").Value = 1.053096140478971
$ws.Range("I16").Value = 1.035721606137383
$ws.Range("J16").Value = 1.042659561163047
$ws.Range("K16").Value = 1.040713607258041
$ws.Range("L16").Value = 1.047817173079473
$ws.Range("M16").Value = 1.056473811899638
$ws.Range("N16").Value = 1.018056102519419

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036785116332067
$ws.Range("D17").Value = 1.037663957923787
$ws.Range("E17").Value = 1.044800996465613
$ws.Range("F17").Value = 1.053524700825711
$ws.Range("I17").Value = 1.035797705427928
$ws.Range("J17").Value = 1.042961449911405
$ws.Range("K17").Value = 1.041024009553957
$ws.Range("L17").Value = 1.048136449360564
$ws.Range("M17").Value = 1.05683058369392
$ws.Range("N17").Value = 1.018158595335275

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037041664621035
$ws.Range("D18").Value = 1.037887566728577
$ws.Range("E18").Value = 1.045029470002198
$ws.Range("F18").Value = 1.053774722235017
$ws.Range("I18").Value = 1.035841918066649
$ws.Range("J18").Value = 1.043137473447782
$ws.Range("K18").Value = 1.041205022833863
$ws.Range("L18").Value = 1.048322647359347
$ws.Range("M18").Value = 1.057038662270929
$ws.Range("N18").Value = 1.01821834153923

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037129140845
$ws.Range("D19").Value = 1.037963817213075
$ws.Range("E19").Value = 1.045107380650786
$ws.Range("F19").Value = 1.053859981465953
$ws.Range("I19").Value = 1.035856963785698
$ws.Range("J19").Value = 1.043197482238182
$ws.Range("K19").Value = 1.041266737056108
$ws.Range("L19").Value = 1.048386130907788
$ws.Range("M19").Value = 1.057109608241695
$ws.Range("N19").Value = 1.018238707324585

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036737926142493
$ws.Range("D20").Value = 1.037622829450908
$ws.Range("E20").Value = 1.044758973834089
$ws.Range("F20").Value = 1.053478715248084
$ws.Range("I20").Value = 1.03578955877336
$ws.Range("J20").Value = 1.042929066630884
$ws.Range("K20").Value = 1.040990710374244
$ws.Range("L20").Value = 1.048102197189056
$ws.Range("M20").Value = 1.056792307566303
$ws.Range("N20").Value = 1.01814760256244

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035466036957898
$ws.Range("D21").Value = 1.0365146365448
$ws.Range("E21").Value = 1.043626758979277
$ws.Range("F21").Value = 1.052239761840338
$ws.Range("I21").Value = 1.035568364270135
$ws.Range("J21").Value = 1.042055679610153
$ws.Range("K21").Value = 1.040092860448905
$ws.Range("L21").Value = 1.047178742672997
$ws.Range("M21").Value = 1.055760492252203
$ws.Range("N21").Value = 1.017850987581035

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03466633488956
$ws.Range("D22").Value = 1.035818170450282
$ws.Range("E22").Value = 1.04291526451255
$ws.Range("F22").Value = 1.051461225351996
$ws.Range("I22").Value = 1.035427716007766
$ws.Range("J22").Value = 1.041505976142029
$ws.Range("K22").Value = 1.039527992172372
$ws.Range("L22").Value = 1.046597853941507
$ws.Range("M22").Value = 1.055111563393399
$ws.Range("N22").Value = 1.017664168407471

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035090273843288
$ws.Range("D23").Value = 1.036187352059185
$ws.Range("E23").Value = 1.043292405821298
$ws.Range("F23").Value = 1.051873900431102
$ws.Range("I23").Value = 1.035502425548502
$ws.Range("J23").Value = 1.04179743852639
$ws.Range("K23").Value = 1.039827473217833
$ws.Range("L23").Value = 1.046905820322479
$ws.Range("M23").Value = 1.055455590554417
$ws.Range("N23").Value = 1.017763235652866

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036759249360166
$ws.Range("D24").Value = 1.037641413534272
$ws.Range("E24").Value = 1.044777961923459
$ws.Range("F24").Value = 1.053499493998303
$ws.Range("I24").Value = 1.035793240435907
$ws.Range("J24").Value = 1.042943699436094
$ws.Range("K24").Value = 1.041005756960929
$ws.Range("L24").Value = 1.048117674364912
$ws.Range("M24").Value = 1.056809602963553
$ws.Range("N24").Value = 1.018152569834865

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038695972742643
$ws.Range("D25").Value = 1.039330079216957
$ws.Range("E25").Value = 1.046503502895315
$ws.Range("F25").Value = 1.055387844239855
$ws.Range("I25").Value = 1.036123858119726
$ws.Range("J25").Value = 1.044271398648893
$ws.Range("K25").Value = 1.042371555635523
$ws.Range("L25").Value = 1.049522769730968
$ws.Range("M25").Value = 1.058380060707272
$ws.Range("N25").Value = 1.018602953538678

